$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Row 18: attempt to add an internal host name without glue (no IP) -> should fail
$lo.ListRows.Add() | Out-Null
$ws.Range("C18").Value = "add"
$ws.Range("D18").Value = "{EMPTY}"
$ws.Range("E18").Value = '{"ns2.epp-16.rst." & $DOMAIN}'
$ws.Range("F18").Value = "{EMPTY}"
$ws.Range("G18").Value = "fail"
$ws.Range("H18").Value = "EPP_UNEXPECTED_COMMAND_SUCCESS"

# Row 19: attempt to add an internal host name with glue (IP) -> should pass
$lo.ListRows.Add() | Out-Null
$ws.Range("C19").Value = "add"
$ws.Range("D19").Value = "{EMPTY}"
$ws.Range("E19").Value = '{"ns2.epp-16.rst." & $DOMAIN}'
$ws.Range("F19").Value = "208.77.190.200"
$ws.Range("G19").Value = "pass"
$ws.Range("H19").Value = "EPP_UNEXPECTED_COMMAND_FAILURE"

# Match the source formatting (bordered style) used for the "hostname" column
# in the other rows of the table that share the same generated-hostname formula.
$ws.Range("E10").Copy()
$ws.Range("E18:E19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H18").Select() | Out-Null
